# Apply the "14 march 1 sit" commit:
#  - Delete Sheet1 (the 5-user locked_out/problem/performance/error/visual fixture)
#  - Keep Sheet3 as the only remaining sheet; it becomes the active sheet/tab
#  - Update Sheet3's data rows: replace the "performance_glitch_user"/"visual_user"
#    rows with a tighter standard_user/problem_user/error_user set
#  - Move the saved selection on the remaining sheet to A5

$wb = $excel.ActiveWorkbook

# Remove the old Sheet1 entirely.
$wb.Worksheets("Sheet1").Delete()

# Work on the sheet that survives.
$ws = $wb.Worksheets("Sheet3")

# Update the three data rows (UserName column) to the new fixture values.
$ws.Range("A2").Value = "standard_user"
$ws.Range("A3").Value = "problem_user"
$ws.Range("A4").Value = "error_user"

# Make Sheet3 the active/selected sheet with the saved selection at A5.
$ws.Activate()
$ws.Range("A5").Select()
